$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 11.02.2022 15:15"

# Row 2 gets refreshed with a new price reading:
#  - B2 becomes the new current price (35.5)
#  - C2 becomes the previous price (old B2 value, 34.9)
#  - D2 becomes the delta, written as a text string "+0.6"
#  - E2 becomes the check timestamp, written as a text string
$ws.Range("B2").Value = 35.5
$ws.Range("C2").Value = 34.9

$ws.Range("D2").Value = "'+0.6"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "2022-02-11 15:15:03"
$ws.Range("E2").Style = "Normal"
